$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range('L2').Value = 2651
$ws.Range('K3').Value = 8183
$ws.Range('L3').Value = 2665
$ws.Range('L4').Value = 713
$ws.Range('L5').Value = 155
$ws.Range('L6').Value = 2394
$ws.Range('K7').Value = 27556
$ws.Range('L7').Value = 8578

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range('L2').Value = 70
$ws.Range('L4').Value = 34
$ws.Range('L5').Value = 33
$ws.Range('L6').Value = 64
$ws.Range('J7').Value = 825
$ws.Range('L7').Value = 282
$ws.Range('L8').Value = 546
$ws.Range('L10').Value = 56
$ws.Range('L11').Value = 150
$ws.Range('L19').Value = 240
$ws.Range('L20').Value = 221
$ws.Range('L23').Value = 85
$ws.Range('L29').Value = 454
$ws.Range('L30').Value = 41
$ws.Range('K33').Value = 1146
$ws.Range('L33').Value = 384
$ws.Range('L36').Value = 117
$ws.Range('L37').Value = 317
$ws.Range('L43').Value = 67
$ws.Range('L47').Value = 64
$ws.Range('L48').Value = 118
$ws.Range('L51').Value = 101
$ws.Range('L52').Value = 172
$ws.Range('L54').Value = 170
$ws.Range('L55').Value = 82
$ws.Range('L60').Value = 52
$ws.Range('J63').Value = 218
$ws.Range('L63').Value = 31
$ws.Range('L64').Value = 56
$ws.Range('L65').Value = 159
$ws.Range('L67').Value = 315
$ws.Range('L69').Value = 23
$ws.Range('L73').Value = 70
$ws.Range('L76').Value = 108
$ws.Range('L78').Value = 113
$ws.Range('J79').Value = 801
$ws.Range('L79').Value = 230
$ws.Range('L83').Value = 201
$ws.Range('L84').Value = 90
$ws.Range('L85').Value = 447
$ws.Range('L87').Value = 25
$ws.Range('L89').Value = 111
$ws.Range('L95').Value = 118
$ws.Range('L97').Value = 77
$ws.Range('K101').Value = 27556
$ws.Range('L101').Value = 8578

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range('L3').Value = 89
$ws.Range('J5').Value = 23
$ws.Range('L6').Value = 80
$ws.Range('J7').Value = 825
$ws.Range('L7').Value = 282

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Range('L2').Value = 55
$ws.Range('L3').Value = 47
$ws.Range('L7').Value = 150

$ws = $wb.Worksheets.Item('Uptown')
$ws.Range('L6').Value = 27
$ws.Range('L7').Value = 111

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range('L2').Value = 132
$ws.Range('L3').Value = 182
$ws.Range('L6').Value = 90
$ws.Range('L7').Value = 447

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range('L2').Value = 58
$ws.Range('L4').Value = 15
$ws.Range('L6').Value = 47
$ws.Range('L7').Value = 172

$ws = $wb.Worksheets.Item('Norwood Park')
$ws.Range('L2').Value = 10
$ws.Range('L3').Value = 5
$ws.Range('L7').Value = 23

$ws = $wb.Worksheets.Item('Austin')
$ws.Range('L2').Value = 155
$ws.Range('L3').Value = 180
$ws.Range('L6').Value = 153
$ws.Range('L7').Value = 546

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range('L3').Value = 86
$ws.Range('L7').Value = 201

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range('L2').Value = 106
$ws.Range('K3').Value = 406
$ws.Range('L3').Value = 120
$ws.Range('L6').Value = 129
$ws.Range('K7').Value = 1146
$ws.Range('L7').Value = 384

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range('L5').Value = 4
$ws.Range('L7').Value = 118

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range('L2').Value = 95
$ws.Range('L7').Value = 317

$ws = $wb.Worksheets.Item('New City')
$ws.Range('L2').Value = 59
$ws.Range('L7').Value = 159

$ws = $wb.Worksheets.Item('Fuller Park')
$ws.Range('L6').Value = 20
$ws.Range('L7').Value = 41

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range('L3').Value = 110
$ws.Range('L6').Value = 77
$ws.Range('L7').Value = 315

$ws = $wb.Worksheets.Item('South Deering')
$ws.Range('L2').Value = 32
$ws.Range('L3').Value = 36
$ws.Range('L6').Value = 20
$ws.Range('L7').Value = 90

$ws = $wb.Worksheets.Item('Loop')
$ws.Range('L6').Value = 88
$ws.Range('L7').Value = 170

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range('L2').Value = 145
$ws.Range('L3').Value = 166
$ws.Range('L6').Value = 118
$ws.Range('L7').Value = 454

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range('L6').Value = 50
$ws.Range('L7').Value = 118

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range('L2').Value = 81
$ws.Range('L7').Value = 240

$ws = $wb.Worksheets.Item('River North')
$ws.Range('L2').Value = 21
$ws.Range('L7').Value = 108

$ws = $wb.Worksheets.Item('Ashburn')
$ws.Range('L4').Value = 6
$ws.Range('L6').Value = 13
$ws.Range('L7').Value = 64

$ws = $wb.Worksheets.Item('Avondale')
$ws.Range('L3').Value = 14
$ws.Range('L7').Value = 56

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range('L3').Value = 32
$ws.Range('L6').Value = 34
$ws.Range('L7').Value = 113

$ws = $wb.Worksheets.Item('Lower West Side')
$ws.Range('L3').Value = 28
$ws.Range('L7').Value = 82

$ws = $wb.Worksheets.Item('Douglas')
$ws.Range('L3').Value = 31
$ws.Range('L7').Value = 85

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range('L3').Value = 85
$ws.Range('J5').Value = 22
$ws.Range('J7').Value = 801
$ws.Range('L7').Value = 230

$ws = $wb.Worksheets.Item('Near South Side')
$ws.Range('L6').Value = 19
$ws.Range('L7').Value = 56

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range('L2').Value = 70
$ws.Range('L7').Value = 221

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range('L6').Value = 28
$ws.Range('L7').Value = 117

$ws = $wb.Worksheets.Item('Kenwood')
$ws.Range('L2').Value = 24
$ws.Range('L7').Value = 64

$ws = $wb.Worksheets.Item('Portage Park')
$ws.Range('L3').Value = 16
$ws.Range('L6').Value = 17
$ws.Range('L7').Value = 70

$ws = $wb.Worksheets.Item('Albany Park')
$ws.Range('L3').Value = 22
$ws.Range('L7').Value = 70

$ws = $wb.Worksheets.Item('West Town')
$ws.Range('L4').Value = 5
$ws.Range('L7').Value = 77

$ws = $wb.Worksheets.Item('Armour Square')
$ws.Range('L2').Value = 9
$ws.Range('L7').Value = 33

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range('L2').Value = 29
$ws.Range('L7').Value = 101

$ws = $wb.Worksheets.Item('Morgan Park')
$ws.Range('L3').Value = 19
$ws.Range('L7').Value = 52

$ws = $wb.Worksheets.Item('Hyde Park')
$ws.Range('L2').Value = 14
$ws.Range('L7').Value = 67

$ws = $wb.Worksheets.Item('Archer Heights')
$ws.Range('L6').Value = 11
$ws.Range('L7').Value = 34

$ws = $wb.Worksheets.Item('Ukrainian Village')
$ws.Range('L6').Value = 12
$ws.Range('L7').Value = 25
